$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @("test1", "test2", "test3", "test4", "test5")
$colB = @("testgg", "testggr", "testggt", "testggy", "testggu")

for ($i = 0; $i -lt $colA.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $colA[$i]
}

for ($i = 0; $i -lt $colB.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}

$ws.Range("D4").Select()
